$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 153
$ws.Range("I12").Value = 135.57143
$ws.Range("K12").Value = 135.57143
$ws.Range("M12").Value = 34.42857000000001
$ws.Range("H33").Value = 306.05264
$ws.Range("I33").Value = 306.05264
$ws.Range("K33").Value = 306.05264
$ws.Range("M33").Value = -77.05264
$ws.Range("H41").Value = 226.06667
$ws.Range("I41").Value = 32.833332
$ws.Range("J41").Value = 999
$ws.Range("K41").Value = 32.833332
$ws.Range("L41").Value = 999
$ws.Range("M41").Value = 407.166668
$ws.Range("N41").Value = -1879
$ws.Range("H43").Value = 22397.15
$ws.Range("J43").Value = 47027
$ws.Range("L43").Value = 47027
$ws.Range("N43").Value = -47165
$ws.Range("H58").Value = 1636.6842
$ws.Range("J58").Value = 2169.5
$ws.Range("L58").Value = 6508.5
$ws.Range("N58").Value = -6808.5
$ws.Range("H63").Value = 60271
$ws.Range("J63").Value = 60271
$ws.Range("L63").Value = 60271
$ws.Range("N63").Value = -61519
$ws.Range("H66").Value = 60271
$ws.Range("J66").Value = 60271
$ws.Range("L66").Value = 180813
$ws.Range("N66").Value = -187053
$ws.Range("H75").Value = 71390
$ws.Range("J75").Value = 71390
$ws.Range("L75").Value = 71390
$ws.Range("N75").Value = -73262
$ws.Range("H78").Value = 71390
$ws.Range("J78").Value = 71390
$ws.Range("L78").Value = 214170
$ws.Range("N78").Value = -223530
$ws.Range("H93").Value = 49497.25
$ws.Range("J93").Value = 49497.25
$ws.Range("L93").Value = 49497.25
$ws.Range("N93").Value = -54489.25
$ws.Range("H95").Value = 35834.715
$ws.Range("J95").Value = 35834.715
$ws.Range("L95").Value = 35834.715
$ws.Range("N95").Value = -41326.715
$ws.Range("H96").Value = 773.25
$ws.Range("J96").Value = 639
$ws.Range("L96").Value = 1917
$ws.Range("N96").Value = -4663
$ws.Range("H100").Value = 2646.25
$ws.Range("I100").Value = 2966.3333
$ws.Range("J100").Value = 2326.1667
$ws.Range("K100").Value = 2966.3333
$ws.Range("L100").Value = 2326.1667
$ws.Range("M100").Value = -2425.3333
$ws.Range("N100").Value = -3408.1667
$ws.Range("H112").Value = 2887.8572
$ws.Range("I112").Value = 1855.3334
$ws.Range("K112").Value = 5566.0002
$ws.Range("M112").Value = -4458.0002
$ws.Range("H132").Value = 446296.53
$ws.Range("I132").Value = 2073.3428
$ws.Range("K132").Value = 6220.028399999999
$ws.Range("M132").Value = -3690.028399999999
$ws.Range("H135").Value = 1565.0834
$ws.Range("I135").Value = 1565.0834
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 14085.7506
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("M135").Value = -11550.7506
$ws.Range("H137").Value = 2772.6562
$ws.Range("J137").Value = 3453.0527
$ws.Range("L137").Value = 10359.1581
$ws.Range("N137").Value = -15459.1581
$ws.Range("H138").Value = 2956.8333
$ws.Range("I138").Value = 3309.1428
$ws.Range("J138").Value = 2904.3618
$ws.Range("K138").Value = 9927.428400000001
$ws.Range("L138").Value = 8713.0854
$ws.Range("M138").Value = -4787.428400000001
$ws.Range("N138").Value = -18993.0854

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23752.623
$ws.Range("I32").Value = 9485.344999999999
$ws.Range("K32").Value = 9485.344999999999
$ws.Range("M32").Value = -9198.344999999999
$ws.Range("H61").Value = 5303
$ws.Range("I61").Value = 6000
$ws.Range("J61").Value = 4954.5
$ws.Range("K61").Value = 6000
$ws.Range("L61").Value = 4954.5
$ws.Range("M61").Value = -5788
$ws.Range("N61").Value = -5378.5
$ws.Range("H74").Value = 998.3333
$ws.Range("I74").Value = 998.3333
$ws.Range("K74").Value = 998.3333
$ws.Range("M74").Value = -124.3333
$ws.Range("H77").Value = 998.3333
$ws.Range("I77").Value = 998.3333
$ws.Range("K77").Value = 4991.6665
$ws.Range("M77").Value = -623.6665000000003
$ws.Range("H110").Value = 2277.6667
$ws.Range("I110").Value = 1854.8823
$ws.Range("K110").Value = 1854.8823
$ws.Range("M110").Value = 190.1177
$ws.Range("H132").Value = 1427.6957
$ws.Range("I132").Value = 919.85
$ws.Range("K132").Value = 2759.55
$ws.Range("M132").Value = -229.5500000000002
$ws.Range("H136").Value = 5303
$ws.Range("I136").Value = 6000
$ws.Range("J136").Value = 4954.5
$ws.Range("K136").Value = 18000
$ws.Range("L136").Value = 14863.5
$ws.Range("M136").Value = -15450
$ws.Range("N136").Value = -19963.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2381.6667
$ws.Range("I20").Value = 2402.4167
$ws.Range("K20").Value = 2402.4167
$ws.Range("M20").Value = -2155.4167
$ws.Range("H22").Value = 684.25
$ws.Range("J22").Value = 1062
$ws.Range("L22").Value = 1062
$ws.Range("N22").Value = -1408
$ws.Range("H24").Value = 1913.6666
$ws.Range("I24").Value = 1913.6666
$ws.Range("K24").Value = 1913.6666
$ws.Range("M24").Value = -1678.6666
$ws.Range("H92").Value = 48950
$ws.Range("J92").Value = 48950
$ws.Range("L92").Value = 48950
$ws.Range("N92").Value = -53942
$ws.Range("H94").Value = 2679.6155
$ws.Range("I94").Value = 1828.8
$ws.Range("K94").Value = 1828.8
$ws.Range("M94").Value = -1377.8
$ws.Range("H105").Value = 2600.1
$ws.Range("I105").Value = 1891.8235
$ws.Range("K105").Value = 1891.8235
$ws.Range("M105").Value = -144.8235
$ws.Range("H134").Value = 3404.3076
$ws.Range("I134").Value = 2841.5908
$ws.Range("K134").Value = 8524.7724
$ws.Range("M134").Value = -5989.7724

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1354.2727
$ws.Range("I16").Value = 1450.8572
$ws.Range("K16").Value = 1450.8572
$ws.Range("M16").Value = -1163.8572
$ws.Range("H31").Value = 5193.476
$ws.Range("I31").Value = 2884.5957
$ws.Range("J31").Value = 11975.8125
$ws.Range("K31").Value = 2884.5957
$ws.Range("L31").Value = 11975.8125
$ws.Range("M31").Value = -2589.5957
$ws.Range("N31").Value = -12565.8125
$ws.Range("H34").Value = 5193.476
$ws.Range("I34").Value = 2884.5957
$ws.Range("J34").Value = 11975.8125
$ws.Range("K34").Value = 2884.5957
$ws.Range("L34").Value = 11975.8125
$ws.Range("M34").Value = -2682.5957
$ws.Range("N34").Value = -12379.8125
$ws.Range("H44").Value = 33665
$ws.Range("J44").Value = 42000
$ws.Range("L44").Value = 42000
$ws.Range("N44").Value = -42884
$ws.Range("H53").Value = 48500
$ws.Range("J53").Value = 48500
$ws.Range("L53").Value = 48500
$ws.Range("N53").Value = -49714
$ws.Range("H58").Value = 4022.4375
$ws.Range("I58").Value = 4022.4375
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 4022.4375
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("M58").Value = -3819.4375
$ws.Range("H74").Value = 55156.5
$ws.Range("J74").Value = 55156.5
$ws.Range("L74").Value = 55156.5
$ws.Range("N74").Value = -56904.5
$ws.Range("H77").Value = 55156.5
$ws.Range("J77").Value = 55156.5
$ws.Range("L77").Value = 165469.5
$ws.Range("N77").Value = -174205.5
$ws.Range("H105").Value = 2465.4827
$ws.Range("J105").Value = 3285.6
$ws.Range("L105").Value = 3285.6
$ws.Range("N105").Value = -6779.6
$ws.Range("H113").Value = 1354.2727
$ws.Range("I113").Value = 1450.8572
$ws.Range("K113").Value = 1450.8572
$ws.Range("M113").Value = 719.1428000000001
$ws.Range("H132").Value = 4486.875
$ws.Range("J132").Value = 4265.6665
$ws.Range("L132").Value = 12796.9995
$ws.Range("N132").Value = -17856.9995
$ws.Range("H134").Value = 3287.75
$ws.Range("I134").Value = 3114.5715
$ws.Range("J134").Value = 4500
$ws.Range("K134").Value = 9343.7145
$ws.Range("L134").Value = 13500
$ws.Range("M134").Value = -6808.7145
$ws.Range("N134").Value = -18570
$ws.Range("H136").Value = 4022.4375
$ws.Range("I136").Value = 4022.4375
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12067.3125
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("M136").Value = -9517.3125
$ws.Range("H141").Value = 469793.6
$ws.Range("J141").Value = 469793.6
$ws.Range("L141").Value = 469793.6
$ws.Range("N141").Value = -480153.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 381
$ws.Range("J2").Value = 496.5
$ws.Range("L2").Value = 2979
$ws.Range("N2").Value = -3205
$ws.Range("H5").Value = 677.6
$ws.Range("I5").Value = 497
$ws.Range("K5").Value = 1491
$ws.Range("M5").Value = -1379
$ws.Range("H7").Value = 44
$ws.Range("I7").Value = 8.666667
$ws.Range("K7").Value = 26.000001
$ws.Range("M7").Value = 85.999999
$ws.Range("H8").Value = 379
$ws.Range("I8").Value = 379
$ws.Range("K8").Value = 1137
$ws.Range("M8").Value = -998
$ws.Range("H32").Value = 250787.25
$ws.Range("J32").Value = 1049.6666
$ws.Range("L32").Value = 3148.9998
$ws.Range("N32").Value = -3714.9998
$ws.Range("H68").Value = 1500
$ws.Range("J68").Value = 2000
$ws.Range("L68").Value = 6000
$ws.Range("N68").Value = -7622
$ws.Range("H71").Value = 1500
$ws.Range("J71").Value = 2000
$ws.Range("L71").Value = 18000
$ws.Range("N71").Value = -26112
$ws.Range("H86").Value = 152327.03
$ws.Range("I86").Value = 556044
$ws.Range("J86").Value = 933.1667
$ws.Range("K86").Value = 1668132
$ws.Range("L86").Value = 2799.5001
$ws.Range("M86").Value = -1666946
$ws.Range("N86").Value = -5171.5001
$ws.Range("H89").Value = 152327.03
$ws.Range("I89").Value = 556044
$ws.Range("J89").Value = 933.1667
$ws.Range("K89").Value = 5004396
$ws.Range("L89").Value = 8398.5003
$ws.Range("M89").Value = -4998468
$ws.Range("N89").Value = -20254.5003
$ws.Range("H114").Value = 793.82355
$ws.Range("I114").Value = 808.5
$ws.Range("J114").Value = 758.6
$ws.Range("K114").Value = 2425.5
$ws.Range("L114").Value = 2275.8
$ws.Range("M114").Value = 828.5
$ws.Range("N114").Value = -8783.799999999999
$ws.Range("H132").Value = 1295.7778
$ws.Range("I132").Value = 1295.7778
$ws.Range("K132").Value = 11662.0002
$ws.Range("M132").Value = -9132.0002
$ws.Range("H135").Value = 677.6
$ws.Range("I135").Value = 497
$ws.Range("K135").Value = 4473
$ws.Range("M135").Value = -1938
$ws.Range("H136").Value = 6281.143
$ws.Range("J136").Value = 2888
$ws.Range("L136").Value = 8664
$ws.Range("N136").Value = -18864

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 2000
$ws.Range("I4").Value = 2000
$ws.Range("K4").Value = 2000
$ws.Range("M4").Value = -1888
$ws.Range("H9").Value = 611.8
$ws.Range("I9").Value = 19.666666
$ws.Range("J9").Value = 1500
$ws.Range("K9").Value = 19.666666
$ws.Range("L9").Value = 1500
$ws.Range("M9").Value = 150.333334
$ws.Range("N9").Value = -1840
$ws.Range("H40").Value = 30018
$ws.Range("J40").Value = 30018
$ws.Range("L40").Value = 30018
$ws.Range("N40").Value = -30320
$ws.Range("H97").Value = 1275.409
$ws.Range("I97").Value = 1223.65
$ws.Range("J97").Value = 1793
$ws.Range("K97").Value = 1223.65
$ws.Range("L97").Value = 1793
$ws.Range("M97").Value = -727.6500000000001
$ws.Range("N97").Value = -2785
$ws.Range("H113").Value = 2409.4285
$ws.Range("I113").Value = 2412.8333
$ws.Range("J113").Value = 2389
$ws.Range("K113").Value = 2412.8333
$ws.Range("L113").Value = 2389
$ws.Range("M113").Value = -242.8332999999998
$ws.Range("N113").Value = -6729
$ws.Range("H132").Value = 2562.5862
$ws.Range("I132").Value = 2678
$ws.Range("K132").Value = 8034
$ws.Range("M132").Value = -5504

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 929.6667
$ws.Range("I22").Value = 929.6667
$ws.Range("K22").Value = 929.6667
$ws.Range("M22").Value = -634.6667
$ws.Range("H27").Value = 929.6667
$ws.Range("I27").Value = 929.6667
$ws.Range("K27").Value = 929.6667
$ws.Range("M27").Value = -822.6667
$ws.Range("H40").Value = 5307.1816
$ws.Range("I40").Value = 5616.364
$ws.Range("K40").Value = 5616.364
$ws.Range("M40").Value = -5480.364
$ws.Range("H46").Value = 1687.1333
$ws.Range("I46").Value = 1968.5
$ws.Range("J46").Value = 1499.5555
$ws.Range("K46").Value = 1968.5
$ws.Range("L46").Value = 1499.5555
$ws.Range("M46").Value = -1780.5
$ws.Range("N46").Value = -1875.5555
$ws.Range("H55").Value = 540.5714
$ws.Range("I55").Value = 704
$ws.Range("K55").Value = 704
$ws.Range("M55").Value = -531
$ws.Range("H61").Value = 1445.5454
$ws.Range("I61").Value = 1445.5454
$ws.Range("K61").Value = 1445.5454
$ws.Range("M61").Value = -1243.5454
$ws.Range("H93").Value = 2835
$ws.Range("I93").Value = 2688.1428
$ws.Range("J93").Value = 3349
$ws.Range("K93").Value = 2688.1428
$ws.Range("L93").Value = 3349
$ws.Range("M93").Value = -1440.1428
$ws.Range("N93").Value = -5845
$ws.Range("H113").Value = 1445.5454
$ws.Range("I113").Value = 1445.5454
$ws.Range("K113").Value = 1445.5454
$ws.Range("M113").Value = 724.4546
$ws.Range("H122").Value = 2906
$ws.Range("I122").Value = 2391.8125
$ws.Range("K122").Value = 7175.4375
$ws.Range("M122").Value = -4725.4375
$ws.Range("H132").Value = 3110.3635
$ws.Range("I132").Value = 3202.5833
$ws.Range("K132").Value = 9607.749899999999
$ws.Range("M132").Value = -7077.749899999999
$ws.Range("H136").Value = 4726.222
$ws.Range("I136").Value = 5443
$ws.Range("J136").Value = 4521.4287
$ws.Range("K136").Value = 16329
$ws.Range("L136").Value = 13564.2861
$ws.Range("M136").Value = -13779
$ws.Range("N136").Value = -18664.2861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 54999.668
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H70").Value = 40000
$ws.Range("J70").Value = 40000
$ws.Range("L70").Value = 40000
$ws.Range("N70").Value = -40630
$ws.Range("H73").Value = 40000
$ws.Range("J73").Value = 40000
$ws.Range("L73").Value = 40000
$ws.Range("N73").Value = -42184
$ws.Range("H96").Value = 1518180.8
$ws.Range("I96").Value = 2427290.5
$ws.Range("J96").Value = 2998
$ws.Range("K96").Value = 2427290.5
$ws.Range("L96").Value = 2998
$ws.Range("M96").Value = -2425917.5
$ws.Range("N96").Value = -5744
$ws.Range("H105").Value = 33038.332
$ws.Range("J105").Value = 33038.332
$ws.Range("L105").Value = 33038.332
$ws.Range("N105").Value = -40026.332
$ws.Range("H107").Value = 960.2727
$ws.Range("I107").Value = 655.5
$ws.Range("K107").Value = 1966.5
$ws.Range("M107").Value = -46.5
$ws.Range("H113").Value = 361.85715
$ws.Range("I113").Value = 351.23077
$ws.Range("K113").Value = 1053.69231
$ws.Range("M113").Value = 1116.30769
$ws.Range("H135").Value = 66474
$ws.Range("J135").Value = 66474
$ws.Range("L135").Value = 66474
$ws.Range("N135").Value = -76614
